$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Point lights" (row 3) and "Point lights in model chain" (row 4) tasks are
# being folded into a single renamed task on row 2 ("Point lights - specular
# lighting"), so delete those two now-redundant rows (comments on B12/B18
# will need to be re-homed onto B10/B16 once the rows above them shift up).
$commentB12 = $ws.Range("B12").Comment.Text()
$commentB18 = $ws.Range("B18").Comment.Text()

$ws.Range("B4").Comment.Delete()
$ws.Range("B12").Comment.Delete()
$ws.Range("B18").Comment.Delete()

$ws.Rows("4").Delete()
$ws.Rows("3").Delete()

# Rename the old "Invesitgate specular issue" task (now row 2) to reflect
# that it now also covers point-light specular handling.
$ws.Range("B2").Value = "Point lights - specular lighting"

# Clear the leftover explicit style on row 2 (it was only ever applying the
# default font, so dropping it collapses back onto the base style).
$ws.Range("A2:C2").ClearFormats()

# Re-add the two surviving comments at their new (shifted-up) locations.
$c10 = $ws.Range("B10").AddComment($commentB12)
$c10.Author = "Jonny"
$c16 = $ws.Range("B16").AddComment($commentB18)
$c16.Author = "Jonny"

$ws.Range("B3").Select()
